# Updated symbol list on Thu Feb  2 20:35:10 UTC 2023 with GitHub Actions
#
# Refreshes the "Price" (column D) and "Volume(1h)" (column E) columns for
# the coin rows whose market data moved since the last scrape. Values are
# plain text (matching the sheet's existing string layout), so each cell is
# forced to Text format before the write and then restored to the default
# "Normal" style so no stray number-format style lingers on the cell (Excel
# would otherwise auto-detect "329.12" / "4.81%" as a number / percentage
# and silently reformat the cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> @(newPrice, newVolume)   (newPrice = $null when the Price cell is unchanged)
$updates = [ordered]@{
    2  = @("329.12",      "4.81%")
    3  = @("40.46",       "8.87%")
    4  = @("5.765",       "12.61%")
    5  = @("0.08122",     "2.46%")
    6  = @("4.593",       "3.96%")
    7  = @("8.778",       "4.47%")
    8  = @("1.975",       "4.65%")
    10 = @("0.9466",      "1.67%")
    11 = @("0.1308",      "5.65%")
    12 = @("0.1993",      "3.97%")
    13 = @("9.019",       "39.26%")
    14 = @("0.09285",     "3.77%")
    15 = @($null,         "4.53%")
    16 = @("0.09637",     "1.19%")
    17 = @($null,         "-5.26%")
    18 = @("0.006327",    "3.67%")
    19 = @("3.359",       "-0.60%")
    20 = @($null,         "2.10%")
    21 = @("0.1422",      "9.56%")
    22 = @("0.2414",      "5.29%")
    24 = @($null,         "5.89%")
    25 = @("0.004360",    "-0.49%")
    26 = @("0.0001093",   "-17.24%")
    27 = @("0.0004000",   "1.34%")
    39 = @("0.02466",     "7.94%")
    40 = @($null,         "3.61%")
    41 = @("0.007465",    "0.12%")
    42 = @("0.1436",      "3.64%")
    43 = @("0.008851",    "4.65%")
    44 = @("0.002054",    "-0.06%")
    45 = @("0.01044",     "33.90%")
    46 = @("0.00006902",  "9.37%")
    47 = @($null,         "1.09%")
    48 = @("0.003514",    "23.63%")
    49 = @($null,         "1.73%")
    50 = @("0.00002105",  "1.09%")
    51 = @("0.0002005",   "1.09%")
}

foreach ($row in $updates.Keys) {
    $price  = $updates[$row][0]
    $volume = $updates[$row][1]

    if ($null -ne $price) {
        $priceCell = $ws.Range("D" + $row)
        $priceCell.NumberFormat = "@"
        $priceCell.Value = $price
        $priceCell.Style = "Normal"
    }

    $volumeCell = $ws.Range("E" + $row)
    $volumeCell.NumberFormat = "@"
    $volumeCell.Value = $volume
    $volumeCell.Style = "Normal"
}
